$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 64 (hunk 0)
$ws.Range("H64").Value = 83336830
$ws.Range("I64").Value = 83336830
$ws.Range("K64").Value = 83336830
$ws.Range("M64").Value = -83336582
# row 67 (hunk 1)
$ws.Range("H67").Value = 83336830
$ws.Range("I67").Value = 83336830
$ws.Range("K67").Value = 83336830
$ws.Range("M67").Value = -83335972
# row 69 (hunk 2)
$ws.Range("H69").Value = 5203.75
$ws.Range("I69").Value = 4900
$ws.Range("K69").Value = 14700
$ws.Range("M69").Value = -13826
# row 72 (hunk 3)
$ws.Range("H72").Value = 5203.75
$ws.Range("I72").Value = 4900
$ws.Range("K72").Value = 44100
$ws.Range("M72").Value = -39732
# row 112 (hunk 4)
$ws.Range("H112").Value = 2612.4285
$ws.Range("J112").Value = 2612.4285
$ws.Range("L112").Value = 7837.2855
$ws.Range("N112").Value = -10053.2855
# row 116 (hunk 5)
$ws.Range("H116").Value = 11121157
$ws.Range("I116").Value = 18530684
$ws.Range("J116").Value = 6867.25
$ws.Range("K116").Value = 18530684
$ws.Range("L116").Value = 6867.25
$ws.Range("M116").Value = -18527242
$ws.Range("N116").Value = -13751.25
# row 129 (hunk 6)
$ws.Range("H129").Value = 55578260
$ws.Range("I129").Value = 9990
$ws.Range("K129").Value = 29970
$ws.Range("M129").Value = -24970
# row 132 (hunk 7)
$ws.Range("H132").Value = 517259.38
$ws.Range("I132").Value = 586623
$ws.Range("K132").Value = 1759869
$ws.Range("M132").Value = -1757339

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 74 (hunk 8)
$ws.Range("H74").Value = 6488.45
$ws.Range("I74").Value = 9081.857
$ws.Range("J74").Value = 5092
$ws.Range("K74").Value = 9081.857
$ws.Range("L74").Value = 5092
$ws.Range("M74").Value = -8207.857
$ws.Range("N74").Value = -6840
# row 77 (hunk 9)
$ws.Range("H77").Value = 6488.45
$ws.Range("I77").Value = 9081.857
$ws.Range("J77").Value = 5092
$ws.Range("K77").Value = 45409.285
$ws.Range("L77").Value = 25460
$ws.Range("M77").Value = -41041.285
$ws.Range("N77").Value = -34196

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 132 (hunk 10)
$ws.Range("H132").Value = 138256.67
$ws.Range("J132").Value = 138256.67
$ws.Range("L132").Value = 138256.67
$ws.Range("N132").Value = -148376.67
# row 134 (hunk 11)
$ws.Range("H134").Value = 1494607.8
$ws.Range("I134").Value = 1769762.9
$ws.Range("J134").Value = 8769.799999999999
$ws.Range("K134").Value = 5309288.699999999
$ws.Range("L134").Value = 26309.4
$ws.Range("M134").Value = -5306753.699999999
$ws.Range("N134").Value = -31379.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 12)
$ws.Range("H31").Value = 5241.2856
$ws.Range("J31").Value = 18765
$ws.Range("L31").Value = 18765
$ws.Range("N31").Value = -19355
# row 34 (hunk 13)
$ws.Range("H34").Value = 5241.2856
$ws.Range("J34").Value = 18765
$ws.Range("L34").Value = 18765
$ws.Range("N34").Value = -19169
# row 58 (hunk 14)
$ws.Range("H58").Value = 23814696
$ws.Range("I58").Value = 32261508
$ws.Range("J58").Value = 10044.454
$ws.Range("K58").Value = 32261508
$ws.Range("L58").Value = 10044.454
$ws.Range("M58").Value = -32261305
$ws.Range("N58").Value = -10450.454
# row 86 (hunk 15)
$ws.Range("H86").Value = 13485.117
$ws.Range("I86").Value = 15929.375
$ws.Range("J86").Value = 11312.444
$ws.Range("K86").Value = 15929.375
$ws.Range("L86").Value = 11312.444
$ws.Range("M86").Value = -14806.375
$ws.Range("N86").Value = -13558.444
# row 89 (hunk 16)
$ws.Range("H89").Value = 13485.117
$ws.Range("I89").Value = 15929.375
$ws.Range("J89").Value = 11312.444
$ws.Range("K89").Value = 79646.875
$ws.Range("L89").Value = 56562.22
$ws.Range("M89").Value = -74030.875
$ws.Range("N89").Value = -67794.22
# row 136 (hunk 17)
$ws.Range("H136").Value = 23814696
$ws.Range("I136").Value = 32261508
$ws.Range("J136").Value = 10044.454
$ws.Range("K136").Value = 96784524
$ws.Range("L136").Value = 30133.362
$ws.Range("M136").Value = -96781974
$ws.Range("N136").Value = -35233.362

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 80 (hunk 18)
$ws.Range("H80").Value = 1900
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# row 83 (hunk 19)
$ws.Range("H83").Value = 1900
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 20 (hunk 20)
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# row 24 (hunk 21)
$ws.Range("H24").Value = 12498.5
$ws.Range("I24").Value = 10000
$ws.Range("J24").Value = 13331.333
$ws.Range("K24").Value = 10000
$ws.Range("L24").Value = 13331.333
$ws.Range("M24").Value = -9827
$ws.Range("N24").Value = -13677.333
# row 70 (hunk 22)
$ws.Range("H70").Value = 16218.0625
$ws.Range("I70").Value = 8724.75
$ws.Range("K70").Value = 8724.75
$ws.Range("M70").Value = -8454.75
# row 73 (hunk 23)
$ws.Range("H73").Value = 16218.0625
$ws.Range("I73").Value = 8724.75
$ws.Range("K73").Value = 8724.75
$ws.Range("M73").Value = -7788.75
# row 110 (hunk 24)
$ws.Range("H110").Value = 64999
$ws.Range("J110").Value = 64999
$ws.Range("L110").Value = 64999
$ws.Range("N110").Value = -73179
# row 133 (hunk 25)
$ws.Range("H133").Value = 149330.67
$ws.Range("J133").Value = 149330.67
$ws.Range("L133").Value = 149330.67
$ws.Range("N133").Value = -159450.67

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 98 (hunk 26)
$ws.Range("H98").Value = 41665.5
$ws.Range("J98").Value = 41665.5
$ws.Range("L98").Value = 41665.5
$ws.Range("N98").Value = -47655.5
# row 136 (hunk 27)
$ws.Range("H136").Value = 9620047
$ws.Range("I136").Value = 31255414
$ws.Range("K136").Value = 93766242
$ws.Range("M136").Value = -93763692

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 62 (hunk 28)
$ws.Range("H62").Value = 17722.857
$ws.Range("I62").Value = 12421.5
$ws.Range("J62").Value = 24791.334
$ws.Range("K62").Value = 12421.5
$ws.Range("L62").Value = 24791.334
$ws.Range("M62").Value = -11797.5
$ws.Range("N62").Value = -26039.334
# row 65 (hunk 29)
$ws.Range("H65").Value = 17722.857
$ws.Range("I65").Value = 12421.5
$ws.Range("J65").Value = 24791.334
$ws.Range("K65").Value = 62107.5
$ws.Range("L65").Value = 123956.67
$ws.Range("M65").Value = -58987.5
$ws.Range("N65").Value = -130196.67
# row 100 (hunk 30)
$ws.Range("H100").Value = 2233.75
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2233.75
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 4467.5
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -5549.5
# row 132 (hunk 31)
$ws.Range("H132").Value = 5256.2446
$ws.Range("I132").Value = 3217.743
$ws.Range("K132").Value = 9653.228999999999
$ws.Range("M132").Value = -7123.228999999999
# row 139 (hunk 32)
$ws.Range("H139").Value = 118999.5
$ws.Range("J139").Value = 118999.5
$ws.Range("L139").Value = 118999.5
$ws.Range("N139").Value = -129279.5
